# Regenerate merged AHB files
#
# 1. Rename the "_old"/"_new" header suffixes to "_FV2310"/"_FV2404".
# 2. Turn the data range into a real Excel Table (ListObject) with an AutoFilter,
#    while preserving the existing header-row formatting (bold/fill/border) so
#    that no incidental dxf / headerRowDxfId ends up in the saved file.
# 3. Freeze the header row (pane split after row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange  = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A1000:U1000")

# Stash the header's current formatting in an out-of-the-way scratch row, then
# strip the header's own formatting so ListObjects.Add doesn't "helpfully"
# capture it into a new header-row dxf.
$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$headerRange.ClearFormats()

# Convert the used range into a worksheet Table with an AutoFilter.
$dataRange = $ws.Range("A1:U62")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Put the original header formatting back, then clean up the scratch row.
$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false
$scratchRange.ClearFormats()

# Freeze the top (header) row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
